$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = "ELT-2A-Circuitos Elétricos 2"
$ws.Range("C3").Value = "-"

# Row 4
$ws.Range("B4").Value = "ELT-2A-Circuitos Elétricos 2"
$ws.Range("C4").Value = "-"
$ws.Range("F4").Value = "ELT-2A-Circuitos Elétricos 2"

# Row 6
$ws.Range("D6").Value = "MCT-2A-Circuitos Elétricos 2"
$ws.Range("E6").Value = "MCT-2A-Programação"

# Row 7
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "-"
